$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new September log entry ("indusind") arrived, so it lands on top of the
# September_Details/September_Date columns (R/S) and every older entry below
# it shifts down by one row. Insert a whole row at 48 to push rows 48-206
# down to 49-207 (this is also what moves the "Broadband" label from A206 to
# A207), then fill in the new entry's two cells.
$ws.Rows("48:48").Insert()

$ws.Cells.Item(48, 18).Value = "indusind"
$ws.Cells.Item(48, 19).Value = "2024-09-24 22:28:01"
